# Control de asistencia - add attendance marks ("x") for Session 3 (column E)
# for the rows that already have attendance marks, matching the
# "mistakes and assistance list" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-based, matching worksheet rows) that receive an "x" in column E
# (Sesion 3) to record attendance.
$rows = @(5, 6, 7, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24)

foreach ($r in $rows) {
    $ws.Range("E$r").Value = "x"
}

# Update the active selection to reflect where the user left off editing.
$ws.Range("E6").Select()
